$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F24").Value = -3.087370019084748
$ws.Range("F28").Value = -3.12828921673713
$ws.Range("F29").Value = -3.133901144269894
$ws.Range("F32").Value = -3.150736926868183
$ws.Range("F34").Value = -3.161960781933712
$ws.Range("F38").Value = -3.181588338794763
$ws.Range("F40").Value = -3.189992040590294
$ws.Range("F43").Value = -3.202597593283588
$ws.Range("F46").Value = -3.194159558595846
$ws.Range("F47").Value = -3.215174766969645
$ws.Range("F48").Value = -3.18715997955034
$ws.Range("F50").Value = -3.19195765643184
$ws.Range("F51").Value = -3.194749430694602
$ws.Range("F52").Value = -3.214354562433405
$ws.Range("F53").Value = -3.1821095451091
$ws.Range("F55").Value = -3.18628301699963
$ws.Range("F56").Value = -3.17085135715136
$ws.Range("F81").Value = -3.311622516534244
$ws.Range("F82").Value = -3.364374510355558
$ws.Range("F83").Value = -3.388847020426648
$ws.Range("F84").Value = -3.380438653880619
$ws.Range("F85").Value = -3.34052531769443
$ws.Range("F86").Value = -3.29593202578861
$ws.Range("F88").Value = -3.169576204225263
$ws.Range("F92").Value = -3.411328268727623
$ws.Range("F93").Value = -3.447413319751786
$ws.Range("F94").Value = -3.455342016166051
$ws.Range("F95").Value = -3.466667505285854
$ws.Range("F96").Value = -3.443359580770089
$ws.Range("F97").Value = -3.410876240705369
$ws.Range("F98").Value = -3.348640728027694
$ws.Range("F99").Value = -3.281491283054691
$ws.Range("F103").Value = -3.483980092655456
$ws.Range("F104").Value = -3.508348672838914
$ws.Range("F105").Value = -3.51562877974288
$ws.Range("F106").Value = -3.51438910820059
$ws.Range("F107").Value = -3.497901544930906
$ws.Range("F108").Value = -3.45642856531793
$ws.Range("F109").Value = -3.40399425219154
$ws.Range("F110").Value = -3.338303325157084
$ws.Range("F111").Value = -3.254875451349846
$ws.Range("F114").Value = -3.53389926994342
$ws.Range("F115").Value = -3.555021073366003
$ws.Range("F116").Value = -3.561333998712417
$ws.Range("F117").Value = -3.557280297221444
$ws.Range("F118").Value = -3.534148750371032
$ws.Range("F119").Value = -3.501774220308993
$ws.Range("F120").Value = -3.450039053641922
$ws.Range("F121").Value = -3.382950290804686
$ws.Range("F122").Value = -3.304295681548265
$ws.Range("F125").Value = -3.571108431704015
$ws.Range("F126").Value = -3.588553562470554
$ws.Range("F127").Value = -3.59692315390874
$ws.Range("F128").Value = -3.590369769702707
$ws.Range("F129").Value = -3.565602464532584
$ws.Range("F130").Value = -3.53030825589324
$ws.Range("F131").Value = -3.479167105616105
$ws.Range("F132").Value = -3.413332562263576
$ws.Range("F133").Value = -3.335415454975907
$ws.Range("F136").Value = -3.597418610859475
$ws.Range("F137").Value = -3.616318582851569
$ws.Range("F138").Value = -3.617965243942504
$ws.Range("F139").Value = -3.613654550752135
$ws.Range("F140").Value = -3.583520445155998
$ws.Range("F141").Value = -3.55586773154011
$ws.Range("F142").Value = -3.505708110746911
$ws.Range("F143").Value = -3.44065786973906
$ws.Range("F144").Value = -3.367724458211546
$ws.Range("F147").Value = -3.620563747279826
$ws.Range("F148").Value = -3.637423272258824
$ws.Range("F149").Value = -3.641359536485421
$ws.Range("F150").Value = -3.63177374549046
$ws.Range("F151").Value = -3.608336112453513
$ws.Range("F152").Value = -3.572365185240318
$ws.Range("F153").Value = -3.525332816345022
$ws.Range("F154").Value = -3.45816437968011
$ws.Range("F155").Value = -3.387075016642637
$ws.Range("F158").Value = -3.641163044890422
$ws.Range("F159").Value = -3.653415322298228
$ws.Range("F160").Value = -3.658026471980685
$ws.Range("F161").Value = -3.646266799484331
$ws.Range("F162").Value = -3.624035075649522
$ws.Range("F163").Value = -3.588909896122948
$ws.Range("F164").Value = -3.531096947512444
$ws.Range("F165").Value = -3.473077165357489
$ws.Range("F166").Value = -3.404731352861566
$ws.Range("F169").Value = -3.654082357887288
$ws.Range("F170").Value = -3.668746245682907
$ws.Range("F171").Value = -3.671195289487139
$ws.Range("F172").Value = -3.664792840949056
$ws.Range("F173").Value = -3.637672124528692
$ws.Range("F174").Value = -3.598506112228528
$ws.Range("F175").Value = -3.546628324161679
$ws.Range("F176").Value = -3.48534452759285
$ws.Range("F177").Value = -3.419107685378437
$ws.Range("F180").Value = -3.669101837013419
$ws.Range("F181").Value = -3.680189309933286
$ws.Range("F182").Value = -3.681372624165619
$ws.Range("F183").Value = -3.670359408408433
$ws.Range("F184").Value = -3.647696516663646
$ws.Range("F185").Value = -3.610745524609876
$ws.Range("F186").Value = -3.556712763049551
$ws.Range("F187").Value = -3.495190291354969
$ws.Range("F188").Value = -3.431007797499183
$ws.Range("F191").Value = -3.678062444949003
$ws.Range("F192").Value = -3.690866314525367
$ws.Range("F193").Value = -3.693537979766068
$ws.Range("F194").Value = -3.68316440080256
$ws.Range("F195").Value = -3.657037778344221
$ws.Range("F196").Value = -3.620742440572214
$ws.Range("F197").Value = -3.568198783152333
$ws.Range("F198").Value = -3.504173091162955
$ws.Range("F199").Value = -3.440684727665591
$ws.Range("F202").Value = -3.688976903388025
$ws.Range("F203").Value = -3.699170357304413
$ws.Range("F204").Value = -3.701061670416574
$ws.Range("F205").Value = -3.68936878145468
$ws.Range("F206").Value = -3.665487094588836
$ws.Range("F207").Value = -3.629229504177753
$ws.Range("F208").Value = -3.577397289823722
$ws.Range("F209").Value = -3.511328444793042
$ws.Range("F210").Value = -3.44914339448307
